$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Insert a new row above row 50, shifting existing rows 50-75 down to 51-76.
$ws.Rows.Item(50).Insert()

# Populate the newly inserted row 50 with the new weekly data entry.
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 45007
$ws.Cells.Item(50, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = 100112031
$ws.Cells.Item(50, 7).Value = "Poroto verde"
$ws.Cells.Item(50, 8).Value = "Magnum"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 100
$ws.Cells.Item(50, 11).Value = 27000
$ws.Cells.Item(50, 12).Value = 28000
$ws.Cells.Item(50, 13).Value = 27500
$ws.Cells.Item(50, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 1100
$ws.Cells.Item(50, 17).Value = 25
$ws.Cells.Item(50, 18).Value = "Hortaliza"
